$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The workbook gained one new weekly price record. Insert a fresh row at
# position 240 (pushing the existing rows 240-283 down to 241-284, which
# also grows the used range from A1:R283 to A1:R284) and populate it with
# the new record's values.
$ws.Rows.Item(240).Insert()

$ws.Cells.Item(240, 1).Value = 9
$ws.Cells.Item(240, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(240, 3).Value = "Metropolitana"
$ws.Cells.Item(240, 4).Value = 44504
$ws.Cells.Item(240, 5).Value = 13
$ws.Cells.Item(240, 6).Value = 100112039
$ws.Cells.Item(240, 7).Value = "Ciboulette"
$ws.Cells.Item(240, 8).Value = "Sin especificar"
$ws.Cells.Item(240, 9).Value = "Primera"
$ws.Cells.Item(240, 10).Value = 206
$ws.Cells.Item(240, 11).Value = 800
$ws.Cells.Item(240, 12).Value = 1000
$ws.Cells.Item(240, 13).Value = 900
$ws.Cells.Item(240, 14).Value = "$/docena de atados"
$ws.Cells.Item(240, 15).Value = "Región Metropolitana"
$ws.Cells.Item(240, 16).Value = 300
$ws.Cells.Item(240, 17).Value = 3
$ws.Cells.Item(240, 18).Value = "Hortaliza"
